$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the affected range as Text so values like "1.01" or "42.834.02"
# are stored as literal strings (matching the source workbook, where these
# "Price"/"Volume" columns are inline strings, not numbers).
$affected = $ws.Range("B2:E51")
$affected.NumberFormat = "@"

$ws.Range('D2').Value = '42.834.02'
$ws.Range('E2').Value = '  +1.24%  '
$ws.Range('D3').Value = '2.284.02'
$ws.Range('E3').Value = '  -0.75%  '
$ws.Range('D4').Value = '1.01'
$ws.Range('E4').Value = '  +0.46%  '
$ws.Range('D5').Value = '313.39'
$ws.Range('E5').Value = '  -0.95%  '
$ws.Range('D6').Value = '104.40'
$ws.Range('E6').Value = '  +1.41%  '
$ws.Range('D7').Value = '0.626'
$ws.Range('E7').Value = '  -0.62%  '
$ws.Range('E8').Value = '  +0.39%  '
$ws.Range('E9').Value = '  +0.04%  '
$ws.Range('D10').Value = '39.33'
$ws.Range('E10').Value = '  -0.74%  '
$ws.Range('D11').Value = '0.0903'
$ws.Range('E11').Value = '  -0.09%  '
$ws.Range('D12').Value = '8.35'
$ws.Range('E12').Value = '  -0.85%  '
$ws.Range('E13').Value = '  +2.62%  '
$ws.Range('D14').Value = '0.991'
$ws.Range('E14').Value = '  +3.27%  '
$ws.Range('D15').Value = '15.21'
$ws.Range('E15').Value = '  -0.23%  '
$ws.Range('D16').Value = '2.631.31'
$ws.Range('E16').Value = '  -0.69%  '
$ws.Range('D17').Value = '2.280.31'
$ws.Range('E17').Value = '  -0.46%  '
$ws.Range('D18').Value = '42.738.19'
$ws.Range('E18').Value = '  +1.02%  '
$ws.Range('D19').Value = '7.40'
$ws.Range('E19').Value = '  -0.64%  '
$ws.Range('D20').Value = '0.0000105'
$ws.Range('E20').Value = '  -0.29%  '
$ws.Range('D21').Value = '13.44'
$ws.Range('E21').Value = '  +17.50%  '
$ws.Range('D22').Value = '73.88'
$ws.Range('E22').Value = '  +0.68%  '
$ws.Range('E23').Value = '  +1.42%  '
$ws.Range('D24').Value = '264.89'
$ws.Range('E24').Value = '  -4.25%  '
$ws.Range('D25').Value = '2.21'
$ws.Range('E25').Value = '  -2.69%  '
$ws.Range('D26').Value = '1.01'
$ws.Range('E26').Value = '  +0.60%  '
$ws.Range('D27').Value = '10.82'
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('D28').Value = '7.12'
$ws.Range('E28').Value = '  +21.17%  '
$ws.Range('E29').Value = '  -0.42%  '
$ws.Range('D30').Value = '22.47'
$ws.Range('E30').Value = '  -1.16%  '
$ws.Range('D31').Value = '37.10'
$ws.Range('E31').Value = '  -1.14%  '
$ws.Range('D32').Value = '166.78'
$ws.Range('E32').Value = '  +0.47%  '
$ws.Range('D33').Value = '0.0872'
$ws.Range('E33').Value = '  -0.29%  '
$ws.Range('E34').Value = '  -3.02%  '
$ws.Range('E35').Value = '  -0.78%  '
$ws.Range('D36').Value = '0.113'
$ws.Range('E36').Value = '  -4.26%  '
$ws.Range('D37').Value = '4.54'
$ws.Range('E37').Value = '  -1.24%  '
$ws.Range('D38').Value = '0.0351'
$ws.Range('E38').Value = '  -4.30%  '
$ws.Range('D39').Value = '3.75'
$ws.Range('E39').Value = '  +0.84%  '
$ws.Range('E40').Value = '  -3.84%  '
$ws.Range('E41').Value = '  +4.28%  '
$ws.Range('B42').Value = 'MultiversX'
$ws.Range('C42').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D42').Value = '70.41'
$ws.Range('E42').Value = '  +0.40%  '
$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').Value = '0.232'
$ws.Range('E43').Value = '  +2.33%  '
$ws.Range('E44').Value = '  +0.46%  '
$ws.Range('D45').Value = '94.11'
$ws.Range('E45').Value = '  -0.58%  '
$ws.Range('D46').Value = '12.09'
$ws.Range('E46').Value = '  +0.17%  '
$ws.Range('D47').Value = '1.740.18'
$ws.Range('E47').Value = '  +9.39%  '
$ws.Range('D48').Value = '112.69'
$ws.Range('E48').Value = '  -0.05%  '
$ws.Range('D49').Value = '79.59'
$ws.Range('E49').Value = '  -1.77%  '
$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D50').Value = '5.20'
$ws.Range('E50').Value = '  -0.49%  '
$ws.Range('B51').Value = 'FraxShare'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D51').Value = '8.73'
$ws.Range('E51').Value = '  -2.78%  '

# Strip the temporary Text format back off so cell styling matches the original
# (un-styled) cells - only the values/content should differ from before.xlsx.
$affected.ClearFormats()
